$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new column N (to the right of the existing M column), inheriting
# the formatting of column M so the new cells pick up the same styles
# (currency/year header/body formats) already used across the table.
$ws.Columns("N").Insert()

# New year column header (2022) and its four data rows.
$ws.Range("N3").Value = 2022
$ws.Range("N4").Value = 6333
$ws.Range("N5").Value = 82675
$ws.Range("N6").Value = 300853

# Match the saved selection/active cell from the edited workbook.
$ws.Range("N2").Select() | Out-Null
